$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 52

# Columns A and D hold zero-padded / date-like text that Excel would
# otherwise auto-coerce into a number/date. Force text entry, then
# restore the default "Normal" style so no extra number-format style
# is left attached to the cell (matches the other data rows, which use
# the default style).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-02-05"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "08:39:01"
$ws.Cells.Item($row, 3).Value = "Wednesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "05"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 125831
$ws.Cells.Item($row, 6).Value = 141823
$ws.Cells.Item($row, 7).Value = 166459
$ws.Cells.Item($row, 8).Value = 157749
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142083
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191006
$ws.Cells.Item($row, 14).Value = 115356
$ws.Cells.Item($row, 15).Value = 44671
$ws.Cells.Item($row, 16).Value = 28219
$ws.Cells.Item($row, 17).Value = 62836
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 38199
$ws.Cells.Item($row, 20).Value = -1
